$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 112, pushing existing rows 112:182 down to 113:183
$ws.Rows.Item(112).Insert()

# Populate the newly inserted row 112 with the new data record
$ws.Range("A112").Value = 3
$ws.Range("B112").Value = "Femacal de La Calera"
$ws.Range("C112").Value = "Coquimbo"
$ws.Range("D112").Value = 44438
$ws.Range("D112").NumberFormat = $ws.Range("D113").NumberFormat
$ws.Range("E112").Value = 5
$ws.Range("F112").Value = 100112009
$ws.Range("G112").Value = "Acelga"
$ws.Range("H112").Value = "Sin especificar"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 310
$ws.Range("K112").Value = 2300
$ws.Range("L112").Value = 2500
$ws.Range("M112").Value = 2397
$ws.Range("N112").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O112").Value = "Provincia de Quillota"
$ws.Range("P112").Value = 400
$ws.Range("Q112").Value = 6
$ws.Range("R112").Value = "Hortaliza"
